# This edit adds two new weekly price records for "Zanahoria" (Carrot) at
# Terminal Hortofrutícola Agro Chillán, inserted right after the header/first
# few rows' context at row 104. All the existing data rows from 104 onward
# are shifted down by two rows (to 106..203), which Excel's row Insert does
# automatically (including carrying the date-column number format s="2").
# The two newly inserted rows (104 and 105) are then populated with their
# own values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 104:105, pushing former rows 104..201 down to 106..203
$ws.Rows("104:105").Insert()

# New row 104
$ws.Range("A104").Value = 7
$ws.Range("B104").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C104").Value = "Ñuble"
$ws.Range("D104").Value = 44512
$ws.Range("E104").Value = 16
$ws.Range("F104").Value = 100114013
$ws.Range("G104").Value = "Zanahoria"
$ws.Range("H104").Value = "Sin especificar"
$ws.Range("I104").Value = "Primera"
$ws.Range("J104").Value = 100
$ws.Range("K104").Value = 8000
$ws.Range("L104").Value = 8500
$ws.Range("M104").Value = 8250
$ws.Range("N104").Value = "$/saco 20 kilos"
$ws.Range("O104").Value = "Provincia de Diguillín"
$ws.Range("P104").Value = 412
$ws.Range("Q104").Value = 20
$ws.Range("R104").Value = "Hortaliza"

# New row 105
$ws.Range("A105").Value = 7
$ws.Range("B105").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C105").Value = "Ñuble"
$ws.Range("D105").Value = 44512
$ws.Range("E105").Value = 16
$ws.Range("F105").Value = 100114013
$ws.Range("G105").Value = "Zanahoria"
$ws.Range("H105").Value = "Sin especificar"
$ws.Range("I105").Value = "Primera"
$ws.Range("J105").Value = 60
$ws.Range("K105").Value = 9000
$ws.Range("L105").Value = 10000
$ws.Range("M105").Value = 9500
$ws.Range("N105").Value = "$/saco 20 kilos"
$ws.Range("O105").Value = "Provincia del Elquí"
$ws.Range("P105").Value = 475
$ws.Range("Q105").Value = 20
$ws.Range("R105").Value = "Hortaliza"
